$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'245.42"
$ws.Cells.Item(2, 5).Value = "'-0.49%"
$ws.Cells.Item(3, 4).Value = "'29.80"
$ws.Cells.Item(3, 5).Value = "'0.26%"
$ws.Cells.Item(4, 4).Value = "'5.205"
$ws.Cells.Item(4, 5).Value = "'0.84%"
$ws.Cells.Item(5, 4).Value = "'0.05746"
$ws.Cells.Item(5, 5).Value = "'0.57%"
$ws.Cells.Item(6, 4).Value = "'6.650"
$ws.Cells.Item(6, 5).Value = "'0.65%"
$ws.Cells.Item(7, 4).Value = "'3.277"
$ws.Cells.Item(7, 5).Value = "'7.21%"
$ws.Cells.Item(8, 4).Value = "'0.8584"
$ws.Cells.Item(8, 5).Value = "'0.02%"
$ws.Cells.Item(9, 4).Value = "'0.8511"
$ws.Cells.Item(9, 5).Value = "'-1.79%"
$ws.Cells.Item(10, 5).Value = "'1.53%"
$ws.Cells.Item(11, 4).Value = "'0.07081"
$ws.Cells.Item(11, 5).Value = "'0.20%"
$ws.Cells.Item(12, 4).Value = "'0.03144"
$ws.Cells.Item(12, 5).Value = "'9.78%"
$ws.Cells.Item(13, 4).Value = "'0.09354"
$ws.Cells.Item(13, 5).Value = "'-0.20%"
$ws.Cells.Item(14, 4).Value = "'0.001530"
$ws.Cells.Item(14, 5).Value = "'0.30%"
$ws.Cells.Item(15, 4).Value = "'0.0005969"
$ws.Cells.Item(15, 5).Value = "'-1.06%"
$ws.Cells.Item(16, 4).Value = "'0.005973"
$ws.Cells.Item(16, 5).Value = "'-3.38%"
$ws.Cells.Item(17, 4).Value = "'3.526"
$ws.Cells.Item(17, 5).Value = "'1.34%"
$ws.Cells.Item(18, 4).Value = "'2.194"
$ws.Cells.Item(18, 5).Value = "'0.89%"
$ws.Cells.Item(19, 4).Value = "'0.3157"
$ws.Cells.Item(19, 5).Value = "'0.34%"
$ws.Cells.Item(20, 5).Value = "'1.78%"
$ws.Cells.Item(21, 5).Value = "'0.57%"
$ws.Cells.Item(22, 4).Value = "'3.491"
$ws.Cells.Item(22, 5).Value = "'0.10%"
$ws.Cells.Item(24, 5).Value = "'-1.15%"
$ws.Cells.Item(25, 4).Value = "'0.001225"
$ws.Cells.Item(25, 5).Value = "'0.26%"
$ws.Cells.Item(26, 4).Value = "'0.004166"
$ws.Cells.Item(26, 5).Value = "'-18.17%"
$ws.Cells.Item(27, 5).Value = "'-0.89%"
$ws.Cells.Item(28, 4).Value = "'0.0001449"
$ws.Cells.Item(28, 5).Value = "'-25.27%"
$ws.Cells.Item(40, 4).Value = "'0.03755"
$ws.Cells.Item(40, 5).Value = "'-0.30%"
$ws.Cells.Item(41, 4).Value = "'0.1072"
$ws.Cells.Item(41, 5).Value = "'0.09%"
$ws.Cells.Item(42, 2).Value = "KickToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(42, 4).Value = "'0.003583"
$ws.Cells.Item(42, 5).Value = "'-37.48%"
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43, 4).Value = "'0.002459"
$ws.Cells.Item(43, 5).Value = "'-5.44%"
$ws.Cells.Item(44, 4).Value = "'0.009924"
$ws.Cells.Item(44, 5).Value = "'1.49%"
$ws.Cells.Item(45, 4).Value = "'0.00005456"
$ws.Cells.Item(45, 5).Value = "'6.70%"
$ws.Cells.Item(46, 5).Value = "'-0.07%"
$ws.Cells.Item(47, 4).Value = "'0.08988"
$ws.Cells.Item(47, 5).Value = "'19.61%"
$ws.Cells.Item(48, 5).Value = "'-18.93%"
$ws.Cells.Item(49, 4).Value = "'0.00002100"
$ws.Cells.Item(49, 5).Value = "'-0.07%"
$ws.Cells.Item(50, 4).Value = "'0.0002000"
$ws.Cells.Item(50, 5).Value = "'-0.07%"
